# Commit: "Tue, Jun 30, 2020  1:05:20 AM"
#
# The underlying OOXML diff shows two things happening to this deck:
#   1) Three tables (on what end up being slides 14, 15 and 16) get their
#      <a:tableStyleId> switched from the deck's single custom table style
#      {1FA163F5-ABA5-49FE-8479-D7649DDB4A4E} to the built-in table style
#      {5AF8F5C3-5217-40E4-98BC-6C3F9E783B7A}.
#   2) The presentation's theme colour scheme changes from the custom
#      "Red Violet" / Integral palette to the stock Office palette (the
#      font scheme and format scheme are untouched - they were already
#      identical between the two theme parts in the package).
#
# Both effects are exactly what happens when a user opens the Design tab
# and clicks the built-in "Office Theme" swatch: PowerPoint re-colours the
# slide master's theme (and, because tables with no explicit style still
# point at the old custom table style, nothing else needs touching there
# except the handful of tables that *were* using a named style id).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-colour the presentation theme: Red Violet -> Office palette.
# ---------------------------------------------------------------------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Index order matches the OOXML <a:clrScheme> child order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# .RGB is a standard COM color long: R + G*256 + B*65536.
$colorScheme.Item(1).RGB  = 0x000000    # dk1      000000
$colorScheme.Item(2).RGB  = 0xFFFFFF    # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 0x6A5444    # dk2      44546A
$colorScheme.Item(4).RGB  = 0xE6E6E7    # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 0xD59B5B    # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 0x317DED    # accent2  ED7D31
$colorScheme.Item(7).RGB  = 0xA5A5A5    # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 0x00C0FF    # accent4  FFC000
$colorScheme.Item(9).RGB  = 0xC47244    # accent5  4472C4
$colorScheme.Item(10).RGB = 0x47AD70    # accent6  70AD47
$colorScheme.Item(11).RGB = 0xC16305    # hlink    0563C1
$colorScheme.Item(12).RGB = 0x724F95    # folHlink 954F72

# ---------------------------------------------------------------------
# 2) Re-apply the (now built-in) table style to the three tables that
#    explicitly carried the old custom style id.
# ---------------------------------------------------------------------
$newTableStyle = "{5AF8F5C3-5217-40E4-98BC-6C3F9E783B7A}"

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq "{1FA163F5-ABA5-49FE-8479-D7649DDB4A4E}") {
                $tbl.ApplyStyle($newTableStyle)
            }
        }
    }
}
